# LL-table.xlsx edit: pseudo-implement <program> and <term> rules
# (see commit message) - updates the LL(1) parsing table on List1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- <def_function> row (row 3): "value" column gained rule 10 ---
$ws.Range("B3").Value = 10

# --- <term_n (row 5): "Ɛ" column rule renumbered 7 -> 9 ---
$ws.Range("O5").Value = 9

# --- <return> row (row 8) becomes <return_value>; its old "return "
#     column entry is dropped and an "expression" column entry is added ---
$ws.Range("A8").Value = "<return_value>"
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 22
$ws.Range("O8").Value = 23

# --- a brand-new <statement_fun> row is inserted at row 9 (what used to
#     be <return_value>'s row slides down and gets renumbered below) ---
$ws.Range("A9").Value = "<statement_fun>"
$ws.Range("D9").Value = 16
$ws.Range("G9").Value = 20
$ws.Range("H9").Value = 15
$ws.Range("I9").Value = 17
$ws.Range("J9").Value = 18
$ws.Range("K9").Value = 19
$ws.Range("O9").Value = 21

# --- <statement> row (row 10): rules renumbered ---
$ws.Range("D10").Value = 25
$ws.Range("H10").Value = 24
$ws.Range("I10").Value = 26
$ws.Range("J10").Value = 27
$ws.Range("K10").Value = 28
$ws.Range("O10").Value = 29

# --- <idwhat> row (row 11): rules renumbered ---
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 31

# --- <assign> row (row 12): rules renumbered ---
$ws.Range("D12").Value = 32
$ws.Range("H12").Value = 33

# --- view state: zoomed to 70%, selection moved to O11 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("O11").Select()
